$aValues = @(
-526.3209500055157,
-524.4479470372619,
-522.5619826905364,
-520.662931608397,
-518.7507575836784,
-516.8253772395135,
-514.8867386148538,
-512.9348661186879,
-510.96971159504,
-508.9913052284678,
-506.9996563762015,
-504.9947869090348,
-502.9767566420569,
-500.9455843324719,
-498.9013113393704,
-496.8440001712912,
-494.7737074712043,
-492.6904892925363,
-490.5944085875326,
-488.485523249671,
-486.3639000992166,
-484.2296142165903,
-482.0827291540471,
-479.923318309864,
-477.7514408685112,
-475.5671951964998,
-473.3706274944925,
-471.1618247999816,
-468.940872925242,
-466.7078104848711,
-464.462736640075,
-462.2057154767515,
-459.9368368608303,
-457.6561677880746,
-455.3637651849463,
-453.0597152692862,
-450.7441062658509,
-448.4169828563089,
-446.0784483024984,
-443.7285680464067,
-441.3673951192543,
-438.995019040299,
-436.6115107136291,
-434.216934671644,
-431.8113915930381,
-429.3949105346652,
-426.9675926260201,
-424.5295056155436,
-422.0806972971806,
-419.6212692187369,
-417.151273553826,
-414.6708024864402,
-412.1798986416407,
-409.6786432714717,
-407.1671067576169,
-404.6453569094813,
-402.1134697778352,
-399.5715032431371,
-397.0195356616037,
-394.4576340006365,
-391.8858629321753,
-389.3042893209173,
-386.7129923121398,
-384.1120294816633,
-381.5014635421405,
-378.8813853470432,
-376.2518483921783,
-373.6129180285574,
-370.9646569795685,
-368.307154618633,
-365.6404627313224,
-362.9646410439397,
-360.2797716173147,
-357.5859128926795,
-354.8831560323578,
-352.1715348921624,
-349.4511317474088,
-346.7220130370501,
-343.9842602169075,
-341.2379157727331,
-338.4830547944313,
-335.7197559666628,
-332.9480846806693,
-330.1681071475258,
-327.3798862233397,
-324.5834879778962,
-321.7789892566498,
-318.9664507762616,
-316.1459347166937,
-313.3175411126151,
-310.4812898374654,
-307.6372920508073,
-304.7855965737814,
-301.9262682298449,
-299.0593834678788,
-296.1850203820026,
-293.3032201123488,
-290.4140890940761,
-287.5176693335586,
-284.6140372226678,
-281.7032703667011,
-278.7854335852587,
-275.8605865910004,
-272.9288117756543,
-269.9901748618988,
-267.0447582494876,
-264.0926103014244,
-261.1338194236249,
-258.1684566317879,
-255.1965819322419,
-252.2182778904401,
-249.2336105298114,
-246.2426507075173,
-243.2454814958002,
-240.2421635749985,
-237.2327790037845,
-234.2173886192221,
-231.1960753389904,
-228.1689049313135,
-225.1359556525492,
-222.0973040578197,
-219.053017207747,
-216.0031779981518,
-212.9478532963885,
-209.8871148678043,
-206.8210523162889,
-203.7497266213126,
-200.6732152873514,
-197.5916029841501,
-194.5049558060952,
-191.4133470276406,
-188.3168605900732,
-185.2155689854326,
-182.1095511846597,
-178.9988856200462,
-175.8836565570712,
-172.7639093763262,
-169.6397563311421,
-166.5112599039525,
-163.3785044221037,
-160.2415637277344,
-157.1005166147034,
-153.9554373999102,
-150.8064161877569,
-147.653521759509,
-144.4968275960122,
-141.3364305742743,
-138.1723992008402,
-135.0048183303786,
-131.83375721047,
-128.6593016616121,
-125.4815424550374,
-122.3005542282172,
-119.1164060595936,
-115.9291840836481,
-112.7389850393786,
-109.5458674599286,
-106.3499282390868,
-103.1512427539577,
-99.94989136854019,
-96.74595524573519,
-93.53952591090258,
-90.3306734055017,
-87.11948939389166,
-83.90604429450177,
-80.69043012719757,
-77.47273335161529,
-74.25302411862549,
-71.03138946779978,
-67.80792701483203,
-64.58269632875601,
-61.35579513112401,
-58.12729945723294,
-54.89730711333954,
-51.66588347739155,
-48.43311547483003,
-45.1990956583507,
-41.9639058674862,
-38.72762217724998,
-35.49033187063916,
-32.25212762785644,
-29.01307509127257,
-25.77328310259912,
-22.5328117746925,
-19.29176483067081,
-16.0502022168778,
-12.80822768719755,
-9.565934266600436,
-6.323376575789477,
-3.080658186022802,
0.1621387999721795,
3.404932694777011,
6.64764139343102,
9.890166270743846,
13.13243680176767,
16.37437626644921,
19.61588648083842,
22.85689095781071,
26.09729111690158,
29.33702208733256,
32.57598511398317,
35.8141095834069,
39.05130279237452,
42.28747721807694,
45.52255116469627,
48.75644877760389,
51.98908800680759,
55.22036555414184,
58.45021599682435,
61.67855783752847,
64.90529887212591,
68.13035187082468,
71.3536410755422,
74.57508495918665,
77.79459689246265,
81.01209007147085,
84.22749412617605,
87.44071160040333,
90.65166686290817,
93.86027896720056,
97.06646117034272,
100.2701385231302,
103.4712300518529,
106.6696438182927,
109.8653102535903,
113.0581378027537,
116.2480460176158,
119.4349650796665,
122.6187928440397,
125.7994709611893,
128.9769025829471,
132.1510135880785,
135.3217232357661,
138.4889630056848,
141.6526360831735,
144.8126687166878,
147.9689804807463,
151.1214921216741,
154.2701296314933,
157.4148022616484,
160.5554475998287,
163.6919750715018,
166.8242995868939,
169.9523756377142,
173.0760908652595,
176.1953773338094,
179.3101617788743,
182.4203652876251,
185.5259150960957,
188.6267209508647,
191.7227269364649,
194.8138378216867,
197.8999863518384,
200.9810944578619,
204.0570892974833,
207.1278926451112,
210.1934266821121,
213.2536304406424,
216.3084143254626,
219.357706582019,
222.4014293796985,
225.4395159321172,
228.4718837673471,
231.4984659070298,
234.5191906639038,
237.5339864325836,
240.5427706694153,
243.5454766767629,
246.5420330488399,
249.5323684509479,
252.5163942853127,
255.4940621730374,
258.4652887660231,
261.4299977916849,
264.3881173961906,
267.3395995102687,
270.2843442681515,
273.222291805249,
276.1533827332928,
279.0775227392694,
281.9946640224617,
284.9047293917728,
287.8076364689289,
290.7033315901436,
293.5917427282365,
296.4727827333097,
299.3464091929567,
302.2125388641067,
305.0711069627077,
307.922041805776,
310.7652669118792,
313.6007282476484,
316.4283544733159,
319.2480568410749,
322.0598051068358,
324.8634982526433,
327.6590735852207,
330.4464801409294,
333.2256316102565,
335.9964707867462,
338.7589257725257,
341.5129330841825,
344.2584219814806,
346.9953184078549,
349.723566078661,
352.4430887492692,
355.1538265970422,
357.8557116656273,
360.5486628908087,
363.2326305648619,
365.907541440833,
368.5733191867066,
371.2299049985633,
373.8772274506714,
376.5152191766381,
379.143820353212,
381.7629531241572,
384.3725586182796,
386.9725488260881,
389.562887656147,
392.14348168136,
394.7142671881351,
397.2751826772408,
399.8261513001669,
402.3671094265976,
404.8979913215935,
407.4187255287925,
409.9292280051077,
412.4294485844933,
414.9193143651427,
417.3987559370753,
419.8676877927811,
422.3260546174901,
424.773782497233,
427.2108005390471,
429.6370494468468,
432.0524319246359,
434.4568947466847,
436.8503685072328,
439.2327546081752,
441.6040199226645,
443.9640628531031,
446.3128150732211,
448.6502083234316,
450.9761801456556,
453.2906459936422,
455.5935294278781,
457.8847616684391,
460.1642660975166,
462.4319596811087,
464.6877853317664,
466.9316593901643,
469.1635083959185,
471.3832617120587,
473.5908491053529,
475.7861775554152,
477.969196724789,
480.1398239906171,
482.2979838085744,
484.4436077164808,
486.5766382787576,
488.6969859144046,
490.8045891470704,
492.8993847100127,
494.9813209202326,
497.0503225650972,
499.1063266070253,
501.1492972948251,
503.1791588015211,
505.1958805930281,
507.1994241214705,
509.1897369028907,
511.166825332004,
513.1306502088964,
515.0812030706855,
517.0185057334626,
518.9425738448315,
520.8534295114689,
522.7511582788068,
524.6358312734219,
526.5075478920714
)

$bValues = @(
-30.540395117831,
-30.28204973086,
-30.02411935436305,
-29.76657914946242,
-29.50941813239587,
-29.25262045067518,
-28.99617589323498,
-28.74008549947851,
-28.48434176895807,
-28.22894837662182,
-27.97390705304349,
-27.71922185605064,
-27.46490224945555,
-27.21095303177219,
-26.95738236595901,
-26.70420145360481,
-26.45142111514548,
-26.1990524264196,
-25.94710759172716,
-25.69559846689049,
-25.44453820935183,
-25.19394113883445,
-24.94382038515704,
-24.6941903850992,
-24.44506407349536,
-24.19645898586833,
-23.94838679597777,
-23.70086391424678,
-23.45390663929405,
-23.20752610101636,
-22.96174020249464,
-22.71656300574676,
-22.4720114059611,
-22.22809983268306,
-21.98484169125899,
-21.7422535250895,
-21.50035204915952,
-21.25914943378112,
-21.01866393663133,
-20.77890986339332,
-20.53990027745272,
-20.30165193554342,
-20.06417966966356,
-19.82749766232021,
-19.59162361993367,
-19.3565681869205,
-19.12234903125713,
-18.88898050781248,
-18.65647505438659,
-18.42485017858561,
-18.19411865486089,
-17.9642969581879,
-17.73539691250798,
-17.50743387972435,
-17.2804221562468,
-17.05437575251338,
-16.82930935775039,
-16.60523598899528,
-16.38217042031008,
-16.16012635527284,
-15.93911724542869,
-15.71915667779089,
-15.50025921369678,
-15.2824375201545,
-15.06570465431207,
-14.8500759282372,
-14.63556348658886,
-14.42218040158252,
-14.20993947358758,
-13.99885555440954,
-13.78894041755495,
-13.58020642634227,
-13.37266761208863,
-13.16633609069029,
-12.96122640091295,
-12.75734865127424,
-12.55471655830144,
-12.35334255731942,
-12.15324010142407,
-11.95441971896247,
-11.7568942552532,
-11.56067670891099,
-11.36577903721661,
-11.17221317926919,
-10.9799907661714,
-10.78912357803981,
-10.59962404078576,
-10.41150339842767,
-10.22477292635503,
-10.03944629791528,
-9.855531842781545,
-9.673043761542232,
-9.491992114416007,
-9.312387926221589,
-9.134242861332183,
-8.957568574072582,
-8.782374369304478,
-8.608673439263953,
-8.43647497538953,
-8.265790159354822,
-8.09663013099032,
-7.929005184455036,
-7.762925244040247,
-7.598401430873139,
-7.435443818375433,
-7.274063303360862,
-7.114268822295969,
-6.956071399989042,
-6.799481003133138,
-6.644506928197046,
-6.491159554795018,
-6.33944829655843,
-6.189382720827098,
-6.040972895647258,
-5.894227669968217,
-5.749156842511297,
-5.605769079843753,
-5.464074002266045,
-5.324080310116329,
-5.185797181211399,
-5.049233596918043,
-4.914398004942427,
-4.781299447260835,
-4.649946149016245,
-4.52034644973466,
-4.392509322648294,
-4.266442424712466,
-4.142153998577811,
-4.019652449949149,
-3.89894532604842,
-3.780040357770247,
-3.662945567117987,
-3.547668417784536,
-3.434216501051367,
-3.322597265289517,
-3.212818240115126,
-3.104885484845454,
-2.998807270285965,
-2.89459001569341,
-2.792240765237494,
-2.691766092703006,
-2.593172614637382,
-2.496466687523224,
-2.401655025690985,
-2.308743569029616,
-2.217738293017675,
-2.128645714717095,
-2.041471408277602,
-1.956221284731328,
-1.872900710326746,
-1.791515371013162,
-1.712070932703426,
-1.634572518511874,
-1.559025007638638,
-1.485433567981147,
-1.413803459004072,
-1.34413894861359,
-1.276445032044251,
-1.210726098685853,
-1.14698651970815,
-1.08523055278647,
-1.025462514763774,
-0.9676861854384811,
-0.9119056108342534,
-0.8581243307010311,
-0.806346122220333,
-0.7565745199581784,
-0.7088126774384023,
-0.6630638797000188,
-0.6193314061202102,
-0.5776178834416896,
-0.5379263193771351,
-0.5002593025870681,
-0.4646195370716271,
-0.4310092266723871,
-0.3994306878458644,
-0.3698861410641143,
-0.3423775592333982,
-0.3169067389530653,
-0.2934754362095759,
-0.2720852943777174,
-0.2527376436616464,
-0.2354339394960618,
-0.2201752221660627,
-0.2069626180762612,
-0.1957969285159606,
-0.186678986657796,
-0.1796094174285372,
-0.1745886373076644,
-0.1716170281304414,
-0.1706947846257378,
-0.1718219682852009,
-0.1749985075959728,
-0.1802241656823753,
-0.1874986023134113,
-0.1968213532459231,
-0.2081917522668846,
-0.2216090556815118,
-0.2370722880839082,
-0.2545804917384089,
-0.2741324037486379,
-0.2957267821927871,
-0.3193621124753179,
-0.3450367624632236,
-0.3727490046171908,
-0.4024970345666645,
-0.4342788674571807,
-0.4680921560841114,
-0.5039347984843232,
-0.5418044078922719,
-0.5816983299453682,
-0.6236138175675308,
-0.6675481373496837,
-0.713498354524944,
-0.7614613248607895,
-0.8114337566029803,
-0.8634124836881369,
-0.9173937729489978,
-0.9733741064718231,
-1.031349751849007,
-1.091316736190748,
-1.153271173502418,
-1.217208940607151,
-1.283125554666157,
-1.351016857126268,
-1.420878100224542,
-1.492704632783756,
-1.566491924367995,
-1.642234563566831,
-1.71992803874997,
-1.799566788767493,
-1.881145673420008,
-1.964659263615886,
-2.050102334573324,
-2.137468773101489,
-2.226752954776032,
-2.31794894893298,
-2.411050731610032,
-2.50605231356522,
-2.602947072631878,
-2.701729132404751,
-2.802391621716055,
-2.904927696662156,
-3.009331707188494,
-3.115595669389279,
-3.223712943503102,
-3.333676604245397,
-3.445479412183033,
-3.55911423110814,
-3.674573160123145,
-3.791849457273514,
-3.910934789087662,
-4.031821761414676,
-4.154502452765769,
-4.278969033742328,
-4.405213344411663,
-4.533227120647553,
-4.663002698213319,
-4.79453109465581,
-4.927803962517174,
-5.062812578746019,
-5.199548601163059,
-5.338002875390487,
-5.478166837964238,
-5.620031611622494,
-5.763588219783969,
-5.908827051402113,
-6.055739144528076,
-6.204315212506856,
-6.354545870393196,
-6.506420747418929,
-6.659931412700359,
-6.815067557274631,
-6.971819115961409,
-7.130176147739178,
-7.290129906457878,
-7.451668748671949,
-7.61478308926894,
-7.779463300544307,
-7.945697689822424,
-8.113477094431687,
-8.282790595441384,
-8.453626852365915,
-8.625976115158664,
-8.799827471275762,
-8.97516898278846,
-9.151991339311451,
-9.330282657163963,
-9.51003209039331,
-9.691228273140467,
-9.873859432317545,
-10.05791514955132,
-10.24338381673675,
-10.43025257089835,
-10.6185123437383,
-10.80814927358824,
-10.99915195517309,
-11.19150972684236,
-11.38520941427814,
-11.58023976793282,
-11.77658843598419,
-11.97424358048488,
-12.17319276608222,
-12.37342322510813,
-12.57492341149455,
-12.77768023718664,
-12.98168175380113,
-13.18691533710347,
-13.39336727615611,
-13.60102618373134,
-13.80987880244433,
-14.01991154113011,
-14.2311121107862,
-14.44346715879169,
-14.65696358085853,
-14.87158882233672,
-15.08732879060398,
-15.30417088051348,
-15.5220999931997,
-15.74110531402619,
-15.96117114647009,
-16.18228422299813,
-16.40443141978318,
-16.62759822354343,
-16.85177120736259,
-17.07693671246278,
-17.30308051645417,
-17.53018719021645,
-17.75824430810582,
-17.9872373496481,
-18.217152057975,
-18.44797261647108,
-18.67968597038731,
-18.9122772899515,
-19.14573198114323,
-19.3800365316913,
-19.61517371173696,
-19.85113051437992,
-20.08789234811206,
-20.32544167432349,
-20.563767446848,
-20.80285175441722,
-21.04267998051483,
-21.28323751058895,
-21.52451033449461,
-21.76648214830672,
-22.009137463229,
-22.25246159673915,
-22.49643904533838,
-22.74105361230855,
-22.98629191932147,
-23.23213774742533,
-23.47857600789712,
-23.72559197031023,
-23.97317101155505,
-24.22129599626894,
-24.46995467320207,
-24.71913123610898,
-24.96881076414058,
-25.21897926681848,
-25.46962419573595,
-25.72072937845752,
-25.97228210576662,
-26.22426967017111,
-26.47668100562902,
-26.72950245754055,
-26.9827221082495,
-27.23633174473943,
-27.49031855995446,
-27.74467579526438,
-27.99939619831866,
-28.25447126028155,
-28.50990063175146,
-28.76567881680103,
-29.02180505102399,
-29.27828346470437,
-29.53511843209753,
-29.79231632309848,
-30.04989339295116,
-30.30786554435505,
-30.56625428074909
)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$startRow = 2
$n = $aValues.Count

for ($i = 0; $i -lt $n; $i++) {
    $r = $startRow + $i
    $ws.Cells.Item($r, 1).Value = $aValues[$i]
    $ws.Cells.Item($r, 2).Value = $bValues[$i]
}

Write-Host "Updated $n rows"
